# Updates the "江西-漫展信息" workbook to add the newly scraped event
# "宜春·第三十七届静卿国风动漫文化汉文化展览会" into both the "展览"
# sheet and the "全部类型" sheet, and refreshes a handful of
# "want to go" counters (column F) that changed between scrapes.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, $row, $col, $text)
    # Force text storage so strings like "2024-10-03" are not
    # auto-converted into date serials by Excel.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

function Insert-Event {
    param($ws, $insertRow, $lastRow, $date, $name, $place, $timeRange, $want, $price, $link, $cover)

    # 1) Push row $insertRow (and everything below it) down by one,
    #    carrying along values/styles exactly as Excel does natively.
    $ws.Rows.Item($insertRow).Insert()

    # 2) Copy the index cell's look (border/bold/centered) from the row
    #    above so the new A-column cell matches the rest of the table.
    $ws.Cells.Item($insertRow - 1, 1).Copy()
    $ws.Cells.Item($insertRow, 1).PasteSpecial(-4122) | Out-Null

    # 3) Populate the new row with the new event's data.
    $ws.Cells.Item($insertRow, 1).Value = $insertRow - 1
    Set-TextCell $ws $insertRow 2 $date
    $ws.Cells.Item($insertRow, 3).Value = $name
    $ws.Cells.Item($insertRow, 4).Value = $place
    $ws.Cells.Item($insertRow, 5).Value = $timeRange
    $ws.Cells.Item($insertRow, 6).Value = $want
    $ws.Cells.Item($insertRow, 7).Value = $price
    $ws.Cells.Item($insertRow, 8).Value = $link
    $ws.Cells.Item($insertRow, 9).Value = $cover

    # 4) The row-number column (A) is a plain sequential literal
    #    (A = row - 1), not a formula, so after shifting everything
    #    below the insertion point down by one row we need to
    #    re-stamp it for every row that moved.
    $r = $insertRow + 1
    while ($r -le $lastRow + 1) {
        $ws.Cells.Item($r, 1).Value = $r - 1
        $r = $r + 1
    }
}

# ---------------------------------------------------------------
# Sheet "展览" (sheet index 1)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# "Want to go" counters that simply grew between scrapes (rows
# unaffected by the later insertion, since they sit above row 24).
$ws1.Cells.Item(12, 6).Value = 213
$ws1.Cells.Item(16, 6).Value = 1848
$ws1.Cells.Item(19, 6).Value = 3725

Insert-Event $ws1 24 37 "2024-10-03" "宜春·第三十七届静卿国风动漫文化汉文化展览会" "宜阳大道19号(交通银行旁) 宜春安缦文华酒店" "2024.10.03 09:00-10.03 17:00" 0 50 "https://show.bilibili.com/platform/detail.html?id=92427" "//i2.hdslb.com/bfs/openplatform/202409/ntPXoGGS1726158288182.jpeg"

# Counters on rows that shifted down by one also grew a little
# between scrapes (new row numbers, after the insert above).
$ws1.Cells.Item(25, 6).Value = 31
$ws1.Cells.Item(26, 6).Value = 2371
$ws1.Cells.Item(37, 6).Value = 1398
$ws1.Cells.Item(38, 6).Value = 131

# ---------------------------------------------------------------
# Sheet "全部类型" (sheet index 4) -- identical table, shifted one
# row further down because of an extra "演出" entry already in it.
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(13, 6).Value = 213
$ws4.Cells.Item(17, 6).Value = 1848
$ws4.Cells.Item(20, 6).Value = 3725

Insert-Event $ws4 25 38 "2024-10-03" "宜春·第三十七届静卿国风动漫文化汉文化展览会" "宜阳大道19号(交通银行旁) 宜春安缦文华酒店" "2024.10.03 09:00-10.03 17:00" 0 50 "https://show.bilibili.com/platform/detail.html?id=92427" "//i2.hdslb.com/bfs/openplatform/202409/ntPXoGGS1726158288182.jpeg"

$ws4.Cells.Item(26, 6).Value = 31
$ws4.Cells.Item(27, 6).Value = 2371
$ws4.Cells.Item(38, 6).Value = 1398
$ws4.Cells.Item(39, 6).Value = 131
